$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.423.12'
$ws.Range("E2").Value = '  +2.23%  '

$ws.Range("D3").Value = '3.268.40'
$ws.Range("E3").Value = '  +1.39%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '398.75'
$ws.Range("E5").Value = '  +0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.97'
$ws.Range("E6").Value = '  -1.63%  '

$ws.Range("E7").Value = '  +4.86%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.620'
$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.25'
$ws.Range("E10").Value = '  -0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0963'
$ws.Range("E11").Value = '  +5.36%  '

$ws.Range("E12").Value = '  +1.19%  '

$ws.Range("D13").Value = '3.780.37'
$ws.Range("E13").Value = '  +1.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.24'
$ws.Range("E14").Value = '  +1.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.96'
$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("D16").Value = '3.254.48'
$ws.Range("E16").Value = '  +0.91%  '

$ws.Range("E17").Value = '  -1.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.12'
$ws.Range("E18").Value = '  +2.79%  '

$ws.Range("D19").Value = '57.237.61'
$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("E20").Value = '  -0.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000108'
$ws.Range("E21").Value = '  +5.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.92'
$ws.Range("E22").Value = '  -0.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '297.18'
$ws.Range("E23").Value = '  -0.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.19'
$ws.Range("E24").Value = '  -1.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.18'
$ws.Range("E25").Value = '  -1.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.16'
$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.38'
$ws.Range("E27").Value = '  +0.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.88'
$ws.Range("E28").Value = '  -3.88%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.43'
$ws.Range("E29").Value = '  -0.48%  '

$ws.Range("E30").Value = '  -1.80%  '

$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("E32").Value = '  +1.25%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.21'
$ws.Range("E33").Value = '  +0.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.34'
$ws.Range("E34").Value = '  +11.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0498'
$ws.Range("E35").Value = '  +1.33%  '

$ws.Range("E36").Value = '  +0.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.44'
$ws.Range("E37").Value = '  +0.35%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.10'
$ws.Range("E38").Value = '  -1.43%  '

$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.05%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.48'
$ws.Range("E40").Value = '  -1.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '138.05'
$ws.Range("E41").Value = '  +2.65%  '

$ws.Range("E42").Value = '  +1.81%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.285'
$ws.Range("E43").Value = '  +0.75%  '

$ws.Range("E44").Value = '  -2.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.85'
$ws.Range("E45").Value = '  -2.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.89'
$ws.Range("E46").Value = '  -2.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.39'
$ws.Range("E47").Value = '  +0.65%  '

$ws.Range("E48").Value = '  +4.43%  '

$ws.Range("D49").Value = '2.157.54'
$ws.Range("E49").Value = '  +1.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.47'
$ws.Range("E50").Value = '  +0.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.94'
$ws.Range("E51").Value = '  -9.63%  '
